# Update end of vignette 1
# 1) Refresh the "datetimeFigureOut" auto-date field (7/15/25 -> 7/23/25)
#    on the slide master and every slide layout's Date Placeholder.
# 2) Slide 3 ("R package research for multicate" two-step approach
#    diagram): reword the "Training data" box from
#    "Multiple RCTs comparing treatment efficacy" to
#    "Multiple studies comparing treatment efficacy", and split the
#    "Step 1" box's second line "Estimate CATE in RCTs" into
#    "Estimate CATE " + "in studies".
# 3) Remove the now-redundant slide title textbox ("TextBox 12":
#    "Two-Step Approach of multicate()").

$p = $ppt.ActivePresentation

# --- 1) Date placeholder refresh -------------------------------------
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "7/23/25"
    }
}

for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $cl = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $shp = $cl.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "7/23/25"
        }
    }
}

# --- 2) Slide 3 text edits ---------------------------------------------
$s3 = $p.Slides.Item(3)

# "Training data" box: Multiple RCTs comparing treatment efficacy -> studies
$rect4 = $s3.Shapes.Item("Rectangle 4")
$para = $rect4.TextFrame.TextRange.Paragraphs(2)
$para.Runs(1).Text = "Multiple studies comparing treatment efficacy"

# "Step 1" box: Estimate CATE in RCTs -> "Estimate CATE " + "in studies"
$rect5 = $s3.Shapes.Item("Rectangle 5")
$para2 = $rect5.TextFrame.TextRange.Paragraphs(2)
$run1 = $para2.Runs(1)
$run1.Text = "Estimate CATE "
$null = $run1.InsertAfter("in studies")

# --- 3) Remove the old slide title textbox -----------------------------
$s3.Shapes.Item("TextBox 12").Delete() | Out-Null
